$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Phase 1 Pre CPP")

# Update the remaining row's values (oxide loading adjustment)
$ws.Range("C2").Value = 0
$ws.Range("E2").Value = 261.6063956874585
$ws.Range("F2").Value = 0.01215375414343174

# Remove rows 3-14, which are no longer part of the data set
$ws.Range("A3:F14").EntireRow.Delete()
